# Weekly update: insert two new fruit/vegetable price survey rows into the
# Cilantro sheet, shifting the existing rows down.
#
# Row 15 (new) is inserted before the old row 15, pushing the old rows
# 15-39 down to 16-40.
# Row 30 (new, in the post-first-insert numbering) is inserted before what
# is now the old row 29 (shifted to 30), pushing it (and everything after)
# down by one more row, to 31-41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row at position 15 -----------------------------
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44427
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = 100112040
$ws.Cells.Item(15, 7).Value = "Cilantro"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 1300
$ws.Cells.Item(15, 12).Value = 1500
$ws.Cells.Item(15, 13).Value = 1400
$ws.Cells.Item(15, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 700
$ws.Cells.Item(15, 17).Value = 2
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# --- Insert the second new row at position 30 -----------------------------
$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(30, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(30, 4).Value = 44417
$ws.Cells.Item(30, 5).Value = 15
$ws.Cells.Item(30, 6).Value = 100112040
$ws.Cells.Item(30, 7).Value = "Cilantro"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 250
$ws.Cells.Item(30, 11).Value = 1000
$ws.Cells.Item(30, 12).Value = 1200
$ws.Cells.Item(30, 13).Value = 1100
$ws.Cells.Item(30, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(30, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 16).Value = 550
$ws.Cells.Item(30, 17).Value = 2
$ws.Cells.Item(30, 18).Value = "Hortaliza"
